# B6-PowerPoint.pptx — commit "Thu, Apr 23, 2020  7:04:54 AM"
#
# 1) Three tables (on slides 14, 15, 16) get their table style switched
#    from the custom "Table_0" style ({3F71EFD1-6FC6-4083-B9AD-BBECC9974C90},
#    also the tableStyles.xml default) to the built-in "No Style, No Grid"
#    style ({183E5761-0AB5-40F6-B1A3-EEB68BCF185A}).
#
# 2) The two theme parts (theme1.xml, used by the slide master, and
#    theme2.xml, used by the notes master) swap contents: theme1 becomes
#    the stock "Office" palette, theme2 becomes the "Integral"/"Red Violet"
#    palette that theme1 used to carry.

$p = $ppt.ActivePresentation

# --- 1) Table style swap --------------------------------------------------

$newTableStyleId = "{183E5761-0AB5-40F6-B1A3-EEB68BCF185A}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2) Theme colour-scheme swap (theme1.xml <-> theme2.xml) -------------
# Only the slide master's theme (theme1.xml) is reachable through the
# PowerPoint object model in this deck (the notes master's theme is not
# independently addressable), so apply the target "Office" palette there.

$colors = $p.Slides.Item(1).ThemeColorScheme

$colors.Item(1).RGB  = 0          # dk1       000000
$colors.Item(2).RGB  = 16777215   # lt1       FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2       44546A
$colors.Item(4).RGB  = 15132391   # lt2       E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1   5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2   ED7D31
$colors.Item(7).RGB  = 10855845   # accent3   A5A5A5
$colors.Item(8).RGB  = 49407      # accent4   FFC000
$colors.Item(9).RGB  = 12874308   # accent5   4472C4
$colors.Item(10).RGB = 4697456    # accent6   70AD47
$colors.Item(11).RGB = 12673797   # hlink     0563C1
$colors.Item(12).RGB = 7491477    # folHlink  954F72
